$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194; this shifts the existing rows
# 194..308 down to 195..309 (matching the dimension change from
# A1:R308 to A1:R309 seen in the diff).
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new record.
$ws.Cells.Item(194, 1).Value  = 4
$ws.Cells.Item(194, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(194, 3).Value  = "Los Lagos"
$ws.Cells.Item(194, 4).Value  = 44529
$ws.Cells.Item(194, 5).Value  = 10
$ws.Cells.Item(194, 6).Value  = 100112006
$ws.Cells.Item(194, 7).Value  = "Repollo"
$ws.Cells.Item(194, 8).Value  = "Crespo record"
$ws.Cells.Item(194, 9).Value  = "Primera"
$ws.Cells.Item(194, 10).Value = 500
$ws.Cells.Item(194, 11).Value = 1200
$ws.Cells.Item(194, 12).Value = 1200
$ws.Cells.Item(194, 13).Value = 1200
$ws.Cells.Item(194, 14).Value = "`$/unidad"
$ws.Cells.Item(194, 15).Value = "Región del Maule"
$ws.Cells.Item(194, 16).Value = 1200
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"
